$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update existing GDP per Capita values in column E (Work Week / Social Spending refresh) ---
$ws.Cells(2, 5).Value = "'700"
$ws.Cells(45, 5).Value = "'1245"
$ws.Cells(82, 5).Value = "'1788"
$ws.Cells(83, 5).Value = "'1808"
$ws.Cells(84, 5).Value = "'1728"
$ws.Cells(85, 5).Value = "'1916"
$ws.Cells(86, 5).Value = "'2099"
$ws.Cells(87, 5).Value = "'1913"
$ws.Cells(88, 5).Value = "'1970"
$ws.Cells(89, 5).Value = "'1978"
$ws.Cells(90, 5).Value = "'1892"
$ws.Cells(91, 5).Value = "'2106"
$ws.Cells(92, 5).Value = "'2197"
$ws.Cells(93, 5).Value = "'2212"
$ws.Cells(94, 5).Value = "'2257"
$ws.Cells(95, 5).Value = "'2270"
$ws.Cells(96, 5).Value = "'2254"
$ws.Cells(97, 5).Value = "'2220"
$ws.Cells(98, 5).Value = "'2158"
$ws.Cells(99, 5).Value = "'2134"
$ws.Cells(100, 5).Value = "'2101"
$ws.Cells(101, 5).Value = "'2112"
$ws.Cells(102, 5).Value = "'2270"
$ws.Cells(103, 5).Value = "'2377"
$ws.Cells(104, 5).Value = "'2235"
$ws.Cells(105, 5).Value = "'2227"
$ws.Cells(106, 5).Value = "'2319"
$ws.Cells(107, 5).Value = "'1988"
$ws.Cells(108, 5).Value = "'1878"
$ws.Cells(109, 5).Value = "'1882"
$ws.Cells(110, 5).Value = "'2008"
$ws.Cells(111, 5).Value = "'1929"
$ws.Cells(112, 5).Value = "'1844"
$ws.Cells(113, 5).Value = "'1820"
$ws.Cells(114, 5).Value = "'1661"
$ws.Cells(115, 5).Value = "'1487"
$ws.Cells(116, 5).Value = "'1530"
$ws.Cells(117, 5).Value = "'1559"
$ws.Cells(118, 5).Value = "'1575"
$ws.Cells(119, 5).Value = "'1605"
$ws.Cells(120, 5).Value = "'1648"
$ws.Cells(121, 5).Value = "'1685"
$ws.Cells(122, 5).Value = "'1693"
$ws.Cells(123, 5).Value = "'1742.62755579766"
$ws.Cells(124, 5).Value = "'1771.60038628535"
$ws.Cells(125, 5).Value = "'1791.11265526599"
$ws.Cells(126, 5).Value = "'1813.6117444642"
$ws.Cells(127, 5).Value = "'1856.57278201958"
$ws.Cells(128, 5).Value = "'1907.93483958593"
$ws.Cells(129, 5).Value = "'1964.91242818097"
$ws.Cells(130, 5).Value = "'2017.98955294768"
$ws.Cells(131, 5).Value = "'2066.74120191922"
$ws.Cells(132, 5).Value = "'2100.29307054518"
$ws.Cells(133, 5).Value = "'2143.10623922565"
$ws.Cells(134, 5).Value = "'2192.0202661421"
$ws.Cells(135, 5).Value = "'2256.81729672897"
$ws.Cells(136, 5).Value = "'2329.5156743128"
$ws.Cells(137, 5).Value = "'2411.28764907302"
$ws.Cells(138, 5).Value = "'2503.95769706586"
$ws.Cells(139, 5).Value = "'2553.16586239898"
$ws.Cells(140, 5).Value = "'2725.78842608382"
$ws.Cells(141, 5).Value = "'2793.87094487804"
$ws.Cells(142, 5).Value = "'2946.03934323595"

# --- Append newly reported years 2011-2016 ---
$ws.Cells(143, 1).Value = 288
$ws.Cells(143, 2).Value = "Ghana"
$ws.Cells(143, 3).Value = "GDP per Capita"
$ws.Cells(143, 4).Value = 2011
$ws.Cells(143, 5).Value = "'3283"

$ws.Cells(144, 1).Value = 288
$ws.Cells(144, 2).Value = "Ghana"
$ws.Cells(144, 3).Value = "GDP per Capita"
$ws.Cells(144, 4).Value = 2012
$ws.Cells(144, 5).Value = "'3508"

$ws.Cells(145, 1).Value = 288
$ws.Cells(145, 2).Value = "Ghana"
$ws.Cells(145, 3).Value = "GDP per Capita"
$ws.Cells(145, 4).Value = 2013
$ws.Cells(145, 5).Value = "'3683"

$ws.Cells(146, 1).Value = 288
$ws.Cells(146, 2).Value = "Ghana"
$ws.Cells(146, 3).Value = "GDP per Capita"
$ws.Cells(146, 4).Value = 2014
$ws.Cells(146, 5).Value = "'3747"

$ws.Cells(147, 1).Value = 288
$ws.Cells(147, 2).Value = "Ghana"
$ws.Cells(147, 3).Value = "GDP per Capita"
$ws.Cells(147, 4).Value = 2015
$ws.Cells(147, 5).Value = "'3809"

$ws.Cells(148, 1).Value = 288
$ws.Cells(148, 2).Value = "Ghana"
$ws.Cells(148, 3).Value = "GDP per Capita"
$ws.Cells(148, 4).Value = 2016
$ws.Cells(148, 5).Value = "'3878"
